$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @{
    2  = 1971
    3  = 704
    4  = 1469
    5  = 1114
    6  = 1166
    7  = 1365
    8  = 377
    9  = 426
    10 = 1188
    11 = 518
    12 = 1013
    14 = 1261
    15 = 697
    16 = 461
    17 = 1601
    18 = 685
    19 = 1542
    20 = 535
    21 = 1483
    22 = 724
    23 = 1446
    24 = 285
    25 = 1677
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 3).Value = $values[$row]
}
